$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.65218286411666981
$ws.Range("E3").Value = 0.94153891426821357
$ws.Range("C4").Value = 0.67948813420904797
$ws.Range("F4").Value = 0.86670710703313969
$ws.Range("D5").Value = 0.74452846691618424
$ws.Range("AV6").Value = 0.72639468444361999
$ws.Range("E7").Value = 0.76626627416271054
$ws.Range("F7").Value = 0.96522711567920805
$ws.Range("H7").Value = 0.97598866651325755
$ws.Range("I7").Value = 0.77213249032962783
$ws.Range("F8").Value = 0.79559211456229106
$ws.Range("H9").Value = 0.99602637499204816
$ws.Range("AA9").Value = 0.80266758663149929
$ws.Range("BO9").Value = 0.80836831617033267
$ws.Range("H10").Value = 0.99757209894392407
$ws.Range("T10").Value = 0.91794318245369011
$ws.Range("J11").Value = 0.89079100970386826
$ws.Range("L11").Value = 0.61198154840837227
$ws.Range("K13").Value = 0.75844298347525285
$ws.Range("L13").Value = 0.98727694315090664
$ws.Range("L14").Value = 0.6597614681261057
$ws.Range("M14").Value = 0.92687106685050114
$ws.Range("O14").Value = 0.79420746299656619
$ws.Range("P14").Value = 0.88970390422441881
$ws.Range("AH14").Value = 0.83068972733752511
$ws.Range("M15").Value = 0.90118397705801634
$ws.Range("Q15").Value = 0.95056510409308381
$ws.Range("O16").Value = 0.90458810120846944
$ws.Range("BH17").Value = 0.84169442575037512
$ws.Range("P18").Value = 0.78037466934243538
$ws.Range("Q18").Value = 0.75235562254996269
$ws.Range("S18").Value = 0.91368850146770875
$ws.Range("T18").Value = 0.85035661934897544
$ws.Range("BC19").Value = 0.73744211347952215
$ws.Range("S21").Value = 0.79790725972304433
$ws.Range("T21").Value = 0.99988027029157944
$ws.Range("BN21").Value = 0.79863034015606726
$ws.Range("T22").Value = 0.7986990993805525
$ws.Range("M23").Value = 0.97817421097217272
$ws.Range("U23").Value = 0.69972993651789117
$ws.Range("Y23").Value = 0.80056587177645355
$ws.Range("X25").Value = 0.72021154142757204
$ws.Range("Z25").Value = 0.71551828197401723
$ws.Range("AA25").Value = 0.73941477082326101
$ws.Range("BN25").Value = 0.98010639671468369
$ws.Range("AA26").Value = 0.6062743451515693
$ws.Range("A27").Value = 0.87553703580610542
$ws.Range("K27").Value = 0.70795770037764949
$ws.Range("AB27").Value = 0.73782391773855938
$ws.Range("Z28").Value = 0.7709034105437691
$ws.Range("AR28").Value = 0.89494251920921708
$ws.Range("V29").Value = 0.62936429912896208
$ws.Range("AC30").Value = 0.90308509275029758
$ws.Range("AC31").Value = 0.59896420369273173
$ws.Range("AD31").Value = 0.82266940050154092
$ws.Range("AF31").Value = 0.93314220195845599
$ws.Range("AG31").Value = 0.82523479754490248
$ws.Range("BJ31").Value = 0.84941216674985387
$ws.Range("AD32").Value = 0.93196421841772881
$ws.Range("BI32").Value = 0.89871736226672405
$ws.Range("AF33").Value = 0.8013599903155455
$ws.Range("BM33").Value = 0.76939296529213053
$ws.Range("AI34").Value = 0.9025138390905878
$ws.Range("AG35").Value = 0.72273984741356767
$ws.Range("AK35").Value = 0.92684308758270606
$ws.Range("AI36").Value = 0.75338552343786547
$ws.Range("AJ37").Value = 0.95489172402098843
$ws.Range("AM37").Value = 0.92788932111207145
$ws.Range("A38").Value = 0.70598663557535413
$ws.Range("AJ38").Value = 0.92820601779961942
$ws.Range("AK38").Value = 0.90382806937096516
$ws.Range("AN38").Value = 0.74256150564179202
$ws.Range("AL39").Value = 0.93515081899326979
$ws.Range("AN39").Value = 0.78316331876721135
$ws.Range("AB40").Value = 0.93628236402949705
$ws.Range("AN41").Value = 0.83644449493560713
$ws.Range("AQ41").Value = 0.72244894253550518
$ws.Range("AN42").Value = 0.79071936920357433
$ws.Range("AO42").Value = 0.98248116739122215
$ws.Range("AQ42").Value = 0.73513744472451692
$ws.Range("B43").Value = 0.73655203091819565
$ws.Range("G43").Value = 0.76879396285022783
$ws.Range("AJ43").Value = 0.68447693268695509
$ws.Range("AP44").Value = 0.72544098733427387
$ws.Range("AS44").Value = 0.94564230079484535
$ws.Range("AT44").Value = 0.93462676994978966
$ws.Range("AZ44").Value = 0.76727811582352135
$ws.Range("AT45").Value = 0.83639261037013479
$ws.Range("AU46").Value = 0.93293939136391901
$ws.Range("BH46").Value = 0.83229288650887445
$ws.Range("AS47").Value = 0.8075731000695332
$ws.Range("AV47").Value = 0.90596265881485172
$ws.Range("AW47").Value = 0.97821520386324379
$ws.Range("AT48").Value = 0.92642463656045804
$ws.Range("AX48").Value = 0.77121916438177696
$ws.Range("P50").Value = 0.9011953121516284
$ws.Range("AA50").Value = 0.83451063134168146
$ws.Range("X51").Value = 0.85210191380940126
$ws.Range("AW51").Value = 0.86584200100986197
$ws.Range("BA51").Value = 0.85329607094332072
$ws.Range("AX52").Value = 0.79710338998759167
$ws.Range("BP52").Value = 0.94111707613561824
$ws.Range("N53").Value = 0.86298233547242642
$ws.Range("V54").Value = 0.65251895312777719
$ws.Range("AW54").Value = 0.85135504655070626
$ws.Range("BA54").Value = 0.80414258874156319
$ws.Range("BC54").Value = 0.94111639968962613
$ws.Range("AM55").Value = 0.76767528625735104
$ws.Range("BD55").Value = 0.68039442159056573
$ws.Range("AY56").Value = 0.75168107836192011
$ws.Range("BB56").Value = 0.90428064783950379
$ws.Range("BE56").Value = 0.76045313744162812
$ws.Range("BF56").Value = 0.75793693882245461
$ws.Range("X57").Value = 0.82454998725378692
$ws.Range("AM57").Value = 0.78582677526975797
$ws.Range("BC57").Value = 0.96122345631486117
$ws.Range("AK58").Value = 0.77952565203670199
$ws.Range("BB59").Value = 0.85281967256401625
$ws.Range("BF59").Value = 0.65924145366984743
$ws.Range("BG60").Value = 0.6899579902675762
$ws.Range("BJ61").Value = 0.9635210840052848
$ws.Range("BK61").Value = 0.83748929451610987
$ws.Range("BH62").Value = 0.92092988628000905
$ws.Range("C63").Value = 0.86218405206990945
$ws.Range("BJ64").Value = 0.84220973008058642
$ws.Range("BK64").Value = 0.836338008345899
$ws.Range("BM64").Value = 0.95772382501379116
$ws.Range("BC65").Value = 0.82769990744862243
$ws.Range("BN65").Value = 0.69293696036053853
$ws.Range("E66").Value = 0.8862784697344428
$ws.Range("F66").Value = 0.8332430136424489
$ws.Range("BL66").Value = 0.65066521597981164
$ws.Range("A67").Value = 0.99322681252559497
$ws.Range("B68").Value = 0.8206131237141856
$ws.Range("BM68").Value = 0.99293433991397351
$ws.Range("BO68").Value = 0.83747921075000742
